# Applies the "Download pdf buttons and contents" edit:
#  1. Education table: turn " (Higher Diploma)" into ", Higher Diploma"
#  2. Education table: turn " (Diploma)" into ", Diploma"
#  3. Education table: turn " (Certificate)" into ", Certificate"
#  4. Remove the empty paragraph + manual page break paragraph that used to
#     separate the SKILLS table from the OPEN SOURCE PROJECTS table.
#  5. Remove the now-stale <w:lastRenderedPageBreak/> marker on the
#     "OPEN SOURCE PROJECTS" heading run (it was an artifact of the manual
#     page break removed in step 4).

$d = $word.ActiveDocument

# --- 1/2/3: education rows -------------------------------------------------
# Each of these substrings is unique across the whole document, so a
# document-scoped Find/Replace is unambiguous. (Using a table-Cell-scoped
# Range for Find/Replace destabilizes this runtime's Paragraphs collection,
# so we deliberately search the whole document instead.)
$d.Content.Find.Execute(" (Higher Diploma)", $false, $false, $false, $false, $false, $true, 1, $false, ", Higher Diploma", 2) | Out-Null
$d.Content.Find.Execute(" (Diploma)", $false, $false, $false, $false, $false, $true, 1, $false, ", Diploma", 2) | Out-Null
$d.Content.Find.Execute(" (Certificate)", $false, $false, $false, $false, $false, $true, 1, $false, ", Certificate", 2) | Out-Null

# --- 4: remove the empty paragraph + the page-break paragraph's run -------
# Locate the paragraph that only contains the manual page break: it is the
# paragraph whose end coincides with the start of the "OPEN SOURCE
# PROJECTS" heading text.
$rngBreak = $d.Content
$rngBreak.Find.Execute("OPEN SOURCE PROJECTS", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$breakParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.End -eq $rngBreak.Start) {
        $breakParaIndex = $i
        break
    }
}

# The paragraph right before it is the stray empty paragraph; delete it
# entirely (matching removal of the stand-alone <w:p/> in the diff).
$pEmpty = $d.Paragraphs.Item($breakParaIndex - 1)
$pEmpty.Range.Delete()

# Re-fetch the (now shifted, same index since the earlier paragraph was
# removed) page-break paragraph and strip out everything but the
# paragraph mark itself (i.e. delete the <w:r><w:br .../></w:r>).
$pBreak2 = $d.Paragraphs.Item($breakParaIndex - 1)
$breakContent = $d.Range($pBreak2.Range.Start, $pBreak2.Range.End - 1)
if ($breakContent.Start -ne $breakContent.End) {
    $breakContent.Delete()
}

# --- 5: drop the lastRenderedPageBreak marker on "OPEN SOURCE PROJECTS" ---
$rngHeading = $d.Content
$rngHeading.Find.Execute("OPEN SOURCE PROJECTS", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$headingStart = $rngHeading.Start
$firstChar = $d.Range($headingStart, $headingStart + 1)
$firstCharText = $firstChar.Text
$firstChar.Delete()
$reinsertPoint = $d.Range($headingStart, $headingStart)
$reinsertPoint.InsertBefore($firstCharText)
